$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -8468
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H107").Value = 352
$ws.Range("I107").Value = 363.5
$ws.Range("K107").Value = 363.5
$ws.Range("M107").Value = 1556.5
$ws.Range("H113").Value = 9936.277
$ws.Range("I113").Value = 8001.3335
$ws.Range("J113").Value = 10323.267
$ws.Range("K113").Value = 8001.3335
$ws.Range("L113").Value = 10323.267
$ws.Range("M113").Value = -4747.3335
$ws.Range("N113").Value = -16831.267
$ws.Range("H116").Value = 6428.4287
$ws.Range("I116").Value = 4999.8335
$ws.Range("J116").Value = 15000
$ws.Range("K116").Value = 4999.8335
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -1557.8335
$ws.Range("N116").Value = -21884
$ws.Range("H132").Value = 18989.816
$ws.Range("I132").Value = 8899.8
$ws.Range("K132").Value = 26699.4
$ws.Range("M132").Value = -24169.4
$ws.Range("H137").Value = 19233454
$ws.Range("I137").Value = 23258400
$ws.Range("J137").Value = 3161.6667
$ws.Range("K137").Value = 69775200
$ws.Range("L137").Value = 9485.000100000001
$ws.Range("M137").Value = -69772650
$ws.Range("N137").Value = -14585.0001
$ws.Range("H138").Value = 2424.238
$ws.Range("I138").Value = 655.86487
$ws.Range("J138").Value = 4940.769
$ws.Range("K138").Value = 1967.59461
$ws.Range("L138").Value = 14822.307
$ws.Range("M138").Value = 3172.40539
$ws.Range("N138").Value = -25102.307
$ws.Range("H141").Value = 5778.5835
$ws.Range("I141").Value = 2934.6
$ws.Range("K141").Value = 8803.8
$ws.Range("M141").Value = -3623.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 2679.5
$ws.Range("I16").Value = 289.33334
$ws.Range("J16").Value = 9850
$ws.Range("K16").Value = 289.33334
$ws.Range("L16").Value = 9850
$ws.Range("M16").Value = -2.333340000000021
$ws.Range("N16").Value = -10424
$ws.Range("H32").Value = 5814.8843
$ws.Range("I32").Value = 6020.161
$ws.Range("K32").Value = 6020.161
$ws.Range("M32").Value = -5733.161
$ws.Range("H61").Value = 2299.35
$ws.Range("I61").Value = 1616.2667
$ws.Range("K61").Value = 1616.2667
$ws.Range("M61").Value = -1404.2667
$ws.Range("H74").Value = 1232.9032
$ws.Range("I74").Value = 1096.2941
$ws.Range("K74").Value = 1096.2941
$ws.Range("M74").Value = -222.2941000000001
$ws.Range("H77").Value = 1232.9032
$ws.Range("I77").Value = 1096.2941
$ws.Range("K77").Value = 5481.4705
$ws.Range("M77").Value = -1113.4705
$ws.Range("H92").Value = 49259.6
$ws.Range("J92").Value = 49259.6
$ws.Range("L92").Value = 49259.6
$ws.Range("N92").Value = -54251.6
$ws.Range("H122").Value = 2498.2856
$ws.Range("I122").Value = 2308.4443
$ws.Range("K122").Value = 6925.3329
$ws.Range("M122").Value = -4475.3329
$ws.Range("H132").Value = 2289.6428
$ws.Range("I132").Value = 2311.8462
$ws.Range("J132").Value = 2001
$ws.Range("K132").Value = 6935.5386
$ws.Range("L132").Value = 6003
$ws.Range("M132").Value = -4405.5386
$ws.Range("N132").Value = -11063
$ws.Range("H136").Value = 2299.35
$ws.Range("I136").Value = 1616.2667
$ws.Range("K136").Value = 4848.800099999999
$ws.Range("M136").Value = -2298.800099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4213.6
$ws.Range("I86").Value = 1836
$ws.Range("J86").Value = 5006.1333
$ws.Range("K86").Value = 1836
$ws.Range("L86").Value = 5006.1333
$ws.Range("M86").Value = -713
$ws.Range("N86").Value = -7252.1333
$ws.Range("H89").Value = 4213.6
$ws.Range("I89").Value = 1836
$ws.Range("J89").Value = 5006.1333
$ws.Range("K89").Value = 9180
$ws.Range("L89").Value = 25030.6665
$ws.Range("M89").Value = -3564
$ws.Range("N89").Value = -36262.66650000001
$ws.Range("H134").Value = 3229.484
$ws.Range("I134").Value = 2909.8333
$ws.Range("J134").Value = 4325.4287
$ws.Range("K134").Value = 8729.499899999999
$ws.Range("L134").Value = 12976.2861
$ws.Range("M134").Value = -6194.499899999999
$ws.Range("N134").Value = -18046.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I22").Value = 435.7143
$ws.Range("J22").Value = 945.6
$ws.Range("K22").Value = 435.7143
$ws.Range("L22").Value = 945.6
$ws.Range("M22").Value = -85.71429999999998
$ws.Range("N22").Value = -1645.6
$ws.Range("H31").Value = 2372.8064
$ws.Range("I31").Value = 2232.9614
$ws.Range("K31").Value = 2232.9614
$ws.Range("M31").Value = -1937.9614
$ws.Range("H34").Value = 2372.8064
$ws.Range("I34").Value = 2232.9614
$ws.Range("K34").Value = 2232.9614
$ws.Range("M34").Value = -2030.9614
$ws.Range("H86").Value = 65997.43
$ws.Range("I86").Value = 85001.5
$ws.Range("J86").Value = 40658.668
$ws.Range("K86").Value = 85001.5
$ws.Range("L86").Value = 40658.668
$ws.Range("M86").Value = -83878.5
$ws.Range("N86").Value = -42904.668
$ws.Range("H89").Value = 65997.43
$ws.Range("I89").Value = 85001.5
$ws.Range("J89").Value = 40658.668
$ws.Range("K89").Value = 425007.5
$ws.Range("L89").Value = 203293.34
$ws.Range("M89").Value = -419391.5
$ws.Range("N89").Value = -214525.34
$ws.Range("H132").Value = 2830.077
$ws.Range("I132").Value = 2881.0454
$ws.Range("K132").Value = 8643.1362
$ws.Range("M132").Value = -6113.136200000001
$ws.Range("H134").Value = 2835.3635
$ws.Range("I134").Value = 2216.375
$ws.Range("J134").Value = 4486
$ws.Range("K134").Value = 6649.125
$ws.Range("L134").Value = 13458
$ws.Range("M134").Value = -4114.125
$ws.Range("N134").Value = -18528

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20120460
$ws.Range("I4").Value = 22864052
$ws.Range("J4").Value = 782.6667
$ws.Range("K4").Value = 68592156
$ws.Range("L4").Value = 2348.0001
$ws.Range("M4").Value = -68592044
$ws.Range("N4").Value = -2572.0001
$ws.Range("H26").Value = 91
$ws.Range("I26").Value = 91
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 273
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 15
$ws.Range("N26").Value = $null
$ws.Range("H97").Value = 461.8
$ws.Range("J97").Value = 477.25
$ws.Range("L97").Value = 1431.75
$ws.Range("N97").Value = -2423.75
$ws.Range("H113").Value = 1544.5834
$ws.Range("J113").Value = 1346.3636
$ws.Range("L113").Value = 4039.0908
$ws.Range("N113").Value = -8379.0908
$ws.Range("H117").Value = 1885.6666
$ws.Range("J117").Value = 2516
$ws.Range("L117").Value = 7548
$ws.Range("N117").Value = -14432
$ws.Range("H121").Value = 989.6429
$ws.Range("I121").Value = 302.85715
$ws.Range("K121").Value = 908.5714499999999
$ws.Range("M121").Value = 401.4285500000001
$ws.Range("H131").Value = 3268.8572
$ws.Range("J131").Value = 4553.5454
$ws.Range("L131").Value = 13660.6362
$ws.Range("N131").Value = -23740.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 15000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H70").Value = 66376.9
$ws.Range("I70").Value = 128746.22
$ws.Range("J70").Value = 10244.5
$ws.Range("K70").Value = 128746.22
$ws.Range("L70").Value = 10244.5
$ws.Range("M70").Value = -128476.22
$ws.Range("N70").Value = -10784.5
$ws.Range("H73").Value = 66376.9
$ws.Range("I73").Value = 128746.22
$ws.Range("J73").Value = 10244.5
$ws.Range("K73").Value = 128746.22
$ws.Range("L73").Value = 10244.5
$ws.Range("M73").Value = -127810.22
$ws.Range("N73").Value = -12116.5
$ws.Range("H126").Value = 1250
$ws.Range("I126").Value = 1250
$ws.Range("K126").Value = 3750
$ws.Range("M126").Value = -1280
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H132").Value = 3275.4783
$ws.Range("I132").Value = 3234.5557
$ws.Range("J132").Value = 3422.8
$ws.Range("K132").Value = 9703.667099999999
$ws.Range("L132").Value = 10268.4
$ws.Range("M132").Value = -7173.667099999999
$ws.Range("N132").Value = -15328.4
$ws.Range("H135").Value = 51184.21
$ws.Range("J135").Value = 51184.21
$ws.Range("L135").Value = 51184.21
$ws.Range("N135").Value = -61324.21

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 2105
$ws.Range("I18").Value = 2105
$ws.Range("K18").Value = 2105
$ws.Range("M18").Value = -1933
$ws.Range("H132").Value = 4584.081
$ws.Range("I132").Value = 3209.8845
$ws.Range("J132").Value = 7832.1816
$ws.Range("K132").Value = 9629.6535
$ws.Range("L132").Value = 23496.5448
$ws.Range("M132").Value = -7099.6535
$ws.Range("N132").Value = -28556.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5852146.5
$ws.Range("I81").Value = 4952.778
$ws.Range("J81").Value = 11114621
$ws.Range("K81").Value = 9905.556
$ws.Range("L81").Value = 22229242
$ws.Range("M81").Value = -8844.556
$ws.Range("N81").Value = -22231364
$ws.Range("H84").Value = 5852146.5
$ws.Range("I84").Value = 4952.778
$ws.Range("J84").Value = 11114621
$ws.Range("K84").Value = 49527.78
$ws.Range("L84").Value = 111146210
$ws.Range("M84").Value = -44223.78
$ws.Range("N84").Value = -111156818
$ws.Range("H122").Value = 706406.2
$ws.Range("I122").Value = 4162.5
$ws.Range("K122").Value = 12487.5
$ws.Range("M122").Value = -10037.5
$ws.Range("H132").Value = 2327.2173
$ws.Range("I132").Value = 2296.682
$ws.Range("K132").Value = 6890.045999999999
$ws.Range("M132").Value = -4360.045999999999
